$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted at row 301 (Fecha 45141 = 2023-08-03),
# pushing all the existing rows 301-422 down by one (to 302-423).
$ws.Rows.Item(301).Insert()

$ws.Cells.Item(301, 1).Value = 10
$ws.Cells.Item(301, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(301, 3).Value = "La Araucanía"
$ws.Cells.Item(301, 4).Value = 45141
$ws.Cells.Item(301, 5).Value = 9
$ws.Cells.Item(301, 6).Value = 100112052
$ws.Cells.Item(301, 7).Value = "Albahaca"
$ws.Cells.Item(301, 8).Value = "Sin especificar"
$ws.Cells.Item(301, 9).Value = "Primera"
$ws.Cells.Item(301, 10).Value = 150
$ws.Cells.Item(301, 11).Value = 6000
$ws.Cells.Item(301, 12).Value = 6000
$ws.Cells.Item(301, 13).Value = 6000
$ws.Cells.Item(301, 14).Value = "$/paquete"
$ws.Cells.Item(301, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(301, 16).Value = 6000
$ws.Cells.Item(301, 17).Value = 1
$ws.Cells.Item(301, 18).Value = "Hortaliza"
